$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 196, pushing existing rows 196:210 down to 197:211
$ws.Rows.Item(196).Insert()

# Populate the newly inserted row 196 with the new record's data
$ws.Cells.Item(196, 1).Value = 5
$ws.Cells.Item(196, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(196, 3).Value = "Maule"
$ws.Cells.Item(196, 4).Value = 45008
$ws.Cells.Item(196, 5).Value = 7
$ws.Cells.Item(196, 6).Value = 100112030
$ws.Cells.Item(196, 7).Value = "Poroto granado"
$ws.Cells.Item(196, 8).Value = "Sin especificar"
$ws.Cells.Item(196, 9).Value = "Primera"
$ws.Cells.Item(196, 10).Value = 300
$ws.Cells.Item(196, 11).Value = 28000
$ws.Cells.Item(196, 12).Value = 28000
$ws.Cells.Item(196, 13).Value = 28000
$ws.Cells.Item(196, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(196, 15).Value = "Región del Maule"
$ws.Cells.Item(196, 16).Value = 1120
$ws.Cells.Item(196, 17).Value = 25
$ws.Cells.Item(196, 18).Value = "Hortaliza"
